$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update supplier discount values
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 13

# Move active selection to B10
$ws.Range("B10").Select()
